$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update values: C6 (lon_max test_value) 10 -> 50
$ws.Range("C6").Value = 50

# Update values: C12 (hr_max test_value) 12 -> 24
$ws.Range("C12").Value = 24

# Update the active cell selection from A9 to C12
$ws.Range("C12").Select()
